# Atualizacao automatica: 2025-08-18 21:00:25
#
# Rows 7-11 get rotated: the detection that used to sit on row 11 moves up
# to row 7, and the previous rows 7-10 each shift down by one row.
# Rows 16-17 get a corrected detection image/bounding box.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-like columns (kept as text / inline strings in the workbook).
$textCols = @("A","D","E","F","I","J")
# Numeric columns (Latitude / Longitude).
$numCols  = @("G","H")

# --- Step 1: snapshot current values for rows 7..11 ---
$rowsData = @{}
for ($r = 7; $r -le 11; $r++) {
    $rowVals = @{}
    foreach ($c in $textCols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    foreach ($c in $numCols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $rowsData[$r] = $rowVals
}

# New order: row7 <- old row11, row8 <- old row7, row9 <- old row8,
#            row10 <- old row9, row11 <- old row10
$newOrder = @{
    7  = 11
    8  = 7
    9  = 8
    10 = 9
    11 = 10
}

foreach ($destRow in 7..11) {
    $srcRow  = $newOrder[$destRow]
    $srcVals = $rowsData[$srcRow]

    # Text columns: prefix with a single-quote so Excel keeps the value as
    # text instead of auto-parsing comma-grouped numbers / decimals.
    foreach ($c in $textCols) {
        $ws.Range("$c$destRow").Value = "'" + $srcVals[$c]
    }

    # Numeric columns: assign as real numbers.
    foreach ($c in $numCols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}

# --- Step 2: update D16/I16 and D17/I17 ---
$ws.Range("D16").Value = "'image_20250807111728_ppp0.jpg"
$ws.Range("I16").Value = "'642,530,686,574"

$ws.Range("D17").Value = "'image_20250807111728_ppp0.jpg"
$ws.Range("I17").Value = "'794,481,830,525"
